$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204170823097229
$ws.Range("B1").Value = 2.522683620452881
$ws.Range("C1").Value = 4.300707817077637
$ws.Range("D1").Value = 2.075307369232178
$ws.Range("E1").Value = 1.175889492034912
